$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New columns: J = SMOKERSTATUS, K = RESERVING_RATE ---------------------
# Header row (row 1)
$ws.Range("J1").Value = "SMOKERSTATUS"
$ws.Range("K1").Value = "RESERVING_RATE"

# Data row (row 2) - policy level smoker status ("U" = unknown) and the
# policy level reserving interest rate (0)
$ws.Range("J2").Value = "U"
$ws.Range("K2").Value = 0

# --- Drop the bespoke font/style that used to highlight the PRODUCT cell ---
# (H2 currently carries a custom Consolas/orange style; the updated sheet
# goes back to the plain default formatting for that cell)
$ws.Range("H2").ClearFormats()

# --- Column widths for the (now used) columns H:K --------------------------
# Best-fit-like widths so the new/affected columns comfortably show their
# header text (engine rounds to the nearest 1/6 character unit).
$ws.Columns("H").ColumnWidth = 16.8
$ws.Columns("I").ColumnWidth = 21.65
$ws.Columns("J").ColumnWidth = 14.15
$ws.Columns("K").ColumnWidth = 15.65

# --- Match the recorded selection in the saved workbook --------------------
$ws.Range("K3").Select() | Out-Null
